$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = -0.711074173183429
$ws.Range("F2").Value = 0.0899049022043186
$ws.Range("G2").Value = -7.90918132103004
$ws.Range("H2").Value = 0.00000000000000259087022348122

# Row 3
$ws.Range("E3").Value = 0.40940071901024
$ws.Range("F3").Value = 0.0520643584448633
$ws.Range("G3").Value = 7.86335856695131
$ws.Range("H3").Value = 0.00000000000000373968542593598

# Row 4
$ws.Range("E4").Value = 0.470594597278706
$ws.Range("F4").Value = 0.0512157161898346
$ws.Range("G4").Value = 9.18848026130133
$ws.Range("H4").Value = 0.0000000000000000000398430386407911

# Row 5
$ws.Range("E5").Value = 0.0898804586003665
$ws.Range("F5").Value = 0.0484139752645723
$ws.Range("G5").Value = 1.85649821377378
$ws.Range("H5").Value = 0.0633825773890801

# Row 6
$ws.Range("E6").Value = -0.0511940379698996
$ws.Range("F6").Value = 0.00881145898349069
$ws.Range("G6").Value = -5.80993886095569
$ws.Range("H6").Value = 0.00000000624956612293004

# Row 7
$ws.Range("E7").Value = 0.352152813946055
